$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the authoritative diff: (CellRef, NewValue, ForceTextFormat)
$updates = @(
    @{Cell='D2'; Value='65.760.37'; ForceText=$false},
    @{Cell='E2'; Value='  +0.61%  '; ForceText=$false},
    @{Cell='D3'; Value='3.382.74'; ForceText=$false},
    @{Cell='E3'; Value='  -0.74%  '; ForceText=$false},
    @{Cell='E4'; Value='  +0.01%  '; ForceText=$false},
    @{Cell='D5'; Value='564.13'; ForceText=$true},
    @{Cell='E5'; Value='  +0.44%  '; ForceText=$false},
    @{Cell='D6'; Value='176.16'; ForceText=$true},
    @{Cell='E6'; Value='  +0.10%  '; ForceText=$false},
    @{Cell='D7'; Value='0.631'; ForceText=$true},
    @{Cell='E7'; Value='  +0.32%  '; ForceText=$false},
    @{Cell='D8'; Value='3.380.02'; ForceText=$false},
    @{Cell='E8'; Value='  -0.55%  '; ForceText=$false},
    @{Cell='E9'; Value='  -0.08%  '; ForceText=$false},
    @{Cell='D10'; Value='0.174'; ForceText=$true},
    @{Cell='E10'; Value='  +1.84%  '; ForceText=$false},
    @{Cell='D11'; Value='0.633'; ForceText=$true},
    @{Cell='E11'; Value='  +0.06%  '; ForceText=$false},
    @{Cell='D12'; Value='53.75'; ForceText=$true},
    @{Cell='E12'; Value='  -2.17%  '; ForceText=$false},
    @{Cell='D13'; Value='0.0000278'; ForceText=$true},
    @{Cell='E13'; Value='  -0.75%  '; ForceText=$false},
    @{Cell='D14'; Value='9.21'; ForceText=$true},
    @{Cell='E14'; Value='  +0.42%  '; ForceText=$false},
    @{Cell='D15'; Value='3.927.05'; ForceText=$false},
    @{Cell='E15'; Value='  -0.48%  '; ForceText=$false},
    @{Cell='B16'; Value='TRON'; ForceText=$false},
    @{Cell='C16'; Value='https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'; ForceText=$false},
    @{Cell='D16'; Value='0.120'; ForceText=$true},
    @{Cell='E16'; Value='  +0.20%  '; ForceText=$false},
    @{Cell='B17'; Value='Chainlink'; ForceText=$false},
    @{Cell='C17'; Value='https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'; ForceText=$false},
    @{Cell='D17'; Value='18.20'; ForceText=$true},
    @{Cell='E17'; Value='  -1.22%  '; ForceText=$false},
    @{Cell='D18'; Value='3.385.30'; ForceText=$false},
    @{Cell='E18'; Value='  -0.61%  '; ForceText=$false},
    @{Cell='D19'; Value='65.867.98'; ForceText=$false},
    @{Cell='E19'; Value='  +0.78%  '; ForceText=$false},
    @{Cell='D20'; Value='11.87'; ForceText=$true},
    @{Cell='E20'; Value='  -0.33%  '; ForceText=$false},
    @{Cell='D21'; Value='0.995'; ForceText=$true},
    @{Cell='E21'; Value='  -0.27%  '; ForceText=$false},
    @{Cell='D22'; Value='463.99'; ForceText=$true},
    @{Cell='E22'; Value='  -1.54%  '; ForceText=$false},
    @{Cell='E23'; Value='  -2.21%  '; ForceText=$false},
    @{Cell='D24'; Value='14.45'; ForceText=$true},
    @{Cell='E24'; Value='  +7.30%  '; ForceText=$false},
    @{Cell='D25'; Value='89.65'; ForceText=$true},
    @{Cell='E25'; Value='  +2.62%  '; ForceText=$false},
    @{Cell='D26'; Value='4.10'; ForceText=$true},
    @{Cell='E26'; Value='  -1.21%  '; ForceText=$false},
    @{Cell='D27'; Value='2.92'; ForceText=$true},
    @{Cell='E27'; Value='  -0.04%  '; ForceText=$false},
    @{Cell='D28'; Value='10.62'; ForceText=$true},
    @{Cell='E28'; Value='  -3.07%  '; ForceText=$false},
    @{Cell='D29'; Value='8.70'; ForceText=$true},
    @{Cell='E29'; Value='  -1.93%  '; ForceText=$false},
    @{Cell='D30'; Value='31.09'; ForceText=$true},
    @{Cell='E30'; Value='  -0.74%  '; ForceText=$false},
    @{Cell='E31'; Value='  -3.11%  '; ForceText=$false},
    @{Cell='B32'; Value='Cosmos'; ForceText=$false},
    @{Cell='C32'; Value='https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'; ForceText=$false},
    @{Cell='D32'; Value='11.46'; ForceText=$true},
    @{Cell='E32'; Value='  -1.03%  '; ForceText=$false},
    @{Cell='B33'; Value='Bittensor'; ForceText=$false},
    @{Cell='C33'; Value='https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'; ForceText=$false},
    @{Cell='D33'; Value='580.57'; ForceText=$true},
    @{Cell='E33'; Value='  +0.73%  '; ForceText=$false},
    @{Cell='D34'; Value='62.20'; ForceText=$true},
    @{Cell='E34'; Value='  -0.28%  '; ForceText=$false},
    @{Cell='E35'; Value='  -0.58%  '; ForceText=$false},
    @{Cell='E36'; Value='  -0.06%  '; ForceText=$false},
    @{Cell='D37'; Value='3.59'; ForceText=$true},
    @{Cell='E37'; Value='  +1.46%  '; ForceText=$false},
    @{Cell='E38'; Value='  +1.06%  '; ForceText=$false},
    @{Cell='D39'; Value='35.97'; ForceText=$true},
    @{Cell='E39'; Value='  +0.01%  '; ForceText=$false},
    @{Cell='D40'; Value='0.377'; ForceText=$true},
    @{Cell='E40'; Value='  +0.87%  '; ForceText=$false},
    @{Cell='D41'; Value='0.0₃0745'; ForceText=$false},
    @{Cell='E41'; Value='  -2.30%  '; ForceText=$false},
    @{Cell='D42'; Value='3.100.78'; ForceText=$false},
    @{Cell='E42'; Value='  +0.03%  '; ForceText=$false},
    @{Cell='E43'; Value='  -1.16%  '; ForceText=$false},
    @{Cell='D44'; Value='0.0417'; ForceText=$true},
    @{Cell='E44'; Value='  -0.21%  '; ForceText=$false},
    @{Cell='E45'; Value='  -1.05%  '; ForceText=$false},
    @{Cell='E46'; Value='  -1.77%  '; ForceText=$false},
    @{Cell='D47'; Value='3.15'; ForceText=$true},
    @{Cell='E47'; Value='  -0.89%  '; ForceText=$false},
    @{Cell='D48'; Value='1.00'; ForceText=$true},
    @{Cell='E48'; Value='  +0.06%  '; ForceText=$false},
    @{Cell='D49'; Value='140.87'; ForceText=$true},
    @{Cell='E49'; Value='  +2.46%  '; ForceText=$false},
    @{Cell='B50'; Value='THORChain'; ForceText=$false},
    @{Cell='C50'; Value='https://coinranking.com/coin/ybmU-kKU+thorchain-rune'; ForceText=$false},
    @{Cell='D50'; Value='8.50'; ForceText=$true},
    @{Cell='E50'; Value='  +2.44%  '; ForceText=$false},
    @{Cell='B51'; Value='WEMIXToken'; ForceText=$false},
    @{Cell='C51'; Value='https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'; ForceText=$false},
    @{Cell='D51'; Value='2.55'; ForceText=$true},
    @{Cell='E51'; Value='  -2.08%  '; ForceText=$false}
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.ForceText) {
        $rng.NumberFormat = '@'
        $rng.Value = $u.Value
        $rng.ClearFormats()
    } else {
        $rng.Value = $u.Value
    }
}
